$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.276.06"
$ws.Range("E2").Value = "'  +0.84%  "
$ws.Range("D3").Value = "'3.737.92"
$ws.Range("E3").Value = "'  -0.04%  "
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("D5").Value = "'612.99"
$ws.Range("E5").Value = "'  +4.63%  "
$ws.Range("D6").Value = "'188.78"
$ws.Range("E6").Value = "'  +5.51%  "
$ws.Range("D7").Value = "'3.734.24"
$ws.Range("E7").Value = "'  +0.13%  "
$ws.Range("E8").Value = "'  +0.40%  "
$ws.Range("E9").Value = "'  +0.51%  "
$ws.Range("D10").Value = "'0.724"
$ws.Range("E10").Value = "'  +0.03%  "
$ws.Range("E11").Value = "'  -3.03%  "
$ws.Range("D12").Value = "'57.73"
$ws.Range("E12").Value = "'  +6.19%  "
$ws.Range("E13").Value = "'  -3.14%  "
$ws.Range("D14").Value = "'10.71"
$ws.Range("E14").Value = "'  -1.29%  "
$ws.Range("D15").Value = "'4.326.08"
$ws.Range("E15").Value = "'  +0.16%  "
$ws.Range("D16").Value = "'3.736.99"
$ws.Range("E16").Value = "'  -0.98%  "
$ws.Range("D17").Value = "'13.14"
$ws.Range("E17").Value = "'  -0.61%  "
$ws.Range("D18").Value = "'19.41"
$ws.Range("E18").Value = "'  -1.07%  "
$ws.Range("E19").Value = "'  -0.40%  "
$ws.Range("E20").Value = "'  -0.94%  "
$ws.Range("D21").Value = "'69.046.37"
$ws.Range("E21").Value = "'  +0.92%  "
$ws.Range("D22").Value = "'415.24"
$ws.Range("E22").Value = "'  +0.15%  "
$ws.Range("D23").Value = "'4.64"
$ws.Range("E23").Value = "'  +0.20%  "
$ws.Range("D24").Value = "'89.51"
$ws.Range("E24").Value = "'  -0.61%  "
$ws.Range("D25").Value = "'3.07"
$ws.Range("E25").Value = "'  -1.28%  "
$ws.Range("D26").Value = "'12.93"
$ws.Range("E26").Value = "'  -0.53%  "
$ws.Range("D27").Value = "'11.00"
$ws.Range("E27").Value = "'  +0.04%  "
$ws.Range("D28").Value = "'6.08"
$ws.Range("E28").Value = "'  +2.15%  "
$ws.Range("E29").Value = "'  -0.86%  "
$ws.Range("E30").Value = "'  +0.08%  "
$ws.Range("D31").Value = "'33.39"
$ws.Range("E31").Value = "'  +0.08%  "
$ws.Range("D32").Value = "'7.37"
$ws.Range("E32").Value = "'  -12.93%  "
$ws.Range("D33").Value = "'12.81"
$ws.Range("E33").Value = "'  -0.69%  "
$ws.Range("D34").Value = "'0.123"
$ws.Range("E34").Value = "'  +2.37%  "
$ws.Range("D35").Value = "'45.11"
$ws.Range("E35").Value = "'  -1.24%  "
$ws.Range("D36").Value = "'619.46"
$ws.Range("E36").Value = "'  +2.18%  "
$ws.Range("D37").Value = "'65.82"
$ws.Range("E37").Value = "'  -1.44%  "
$ws.Range("D38").Value = "'0.0₃0853"
$ws.Range("E38").Value = "'  -10.86%  "
$ws.Range("D39").Value = "'0.413"
$ws.Range("E39").Value = "'  +0.73%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "'  -0.01%  "
$ws.Range("E41").Value = "'  +0.11%  "
$ws.Range("E42").Value = "'  +2.83%  "
$ws.Range("E43").Value = "'  -0.60%  "
$ws.Range("E44").Value = "'  -0.51%  "
$ws.Range("D45").Value = "'2.66"
$ws.Range("E46").Value = "'  +3.89%  "
$ws.Range("D47").Value = "'9.25"
$ws.Range("E47").Value = "'  -4.02%  "
$ws.Range("D48").Value = "'2.834.56"
$ws.Range("E48").Value = "'  +2.64%  "
$ws.Range("E49").Value = "'  +4.29%  "
$ws.Range("D50").Value = "'2.73"
$ws.Range("E50").Value = "'  -20.48%  "
$ws.Range("D51").Value = "'3.15"
$ws.Range("E51").Value = "'  -3.38%  "
